$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Split the "Ministry Course Code and Level" column (G) into two columns:
#   G: "Ministry Course Code" (e.g. "ENST")
#   H: "Ministry Course Level" (e.g. 12)
# Insert a new column at H so existing G keeps its code value/format,
# and everything that was to the right of G shifts one column over.
$ws.Columns("H").Insert()

# Headers
$ws.Range("G1").Value = "Ministry Course Code"
$ws.Range("G2").Value = "ENST"
$ws.Range("H1").Value = "Ministry Course Level"
$ws.Range("H2").Value = 12

$ws.Range("G3").Value = "ENST"
$ws.Range("H3").Value = 12

$ws.Range("G4").Value = "ENST"
$ws.Range("H4").Value = 12

# Reflect the user having just finished editing/reviewing the new columns
$ws.Columns("G:H").Select() | Out-Null
